# Insert a new data row at row 527 (pushing the existing row 527 and all
# subsequent rows down by one, through the former last row 605 which
# becomes row 606), and populate the newly inserted row with the new
# weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 527.
$ws.Rows.Item(527).Insert()

# Fill in the values for the newly inserted row 527.
$ws.Cells.Item(527, 1).Value = 9
$ws.Cells.Item(527, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(527, 3).Value = "Metropolitana"
$ws.Cells.Item(527, 4).Value = 45127
$ws.Cells.Item(527, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(527, 5).Value = 13
$ws.Cells.Item(527, 6).Value = 100112052
$ws.Cells.Item(527, 7).Value = "Albahaca"
$ws.Cells.Item(527, 8).Value = "Sin especificar"
$ws.Cells.Item(527, 9).Value = "Primera"
$ws.Cells.Item(527, 10).Value = 250
$ws.Cells.Item(527, 11).Value = 4000
$ws.Cells.Item(527, 12).Value = 4500
$ws.Cells.Item(527, 13).Value = 4250
$ws.Cells.Item(527, 14).Value = "`$/paquete"
$ws.Cells.Item(527, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(527, 16).Value = 4250
$ws.Cells.Item(527, 17).Value = 1
$ws.Cells.Item(527, 18).Value = "Hortaliza"
